# Thanh_ERD.docx: drop the two lead-in slides ("Con người tham gia hệ
# thống:" + its picture, and "Thông tin thu thập được:" + its picture,
# plus the page breaks that separated them) so the document opens
# directly on the big ERD diagram that used to be the third page.

$d = $word.ActiveDocument

# Paragraphs 1-6 are, in document order:
#   1. "Con người tham gia hệ thống:" heading
#   2. Picture 1 (inline drawing)
#   3. page break
#   4. "Thông tin thu thập được:" heading
#   5. Picture 2 (inline drawing)
#   6. page break
# Paragraph 7 (which becomes paragraph 1 afterwards) holds the ERD
# diagram/quiz artwork and must be left untouched.
$firstPara = $d.Paragraphs.Item(1)
$lastPara = $d.Paragraphs.Item(6)

$deleteRange = $d.Range($firstPara.Range.Start, $lastPara.Range.End)
$deleteRange.Delete()
